# Add the missing "2022-Q4" quarterly sheet (inserted right after "总计",
# pushing 2022-Q3 / 2022-Q2 / 2022-Q1 / 2021-Q4 one slot later) and record
# its summary row at the top of the "总计" sheet.

$wb = $excel.ActiveWorkbook

# A cell that already carries the workbook's "bold, centered, boxed" header /
# index-column style, used below to stamp the same look onto new cells via
# copy/paste-special so we reuse the existing style instead of inventing one.
$total = $wb.Worksheets.Item(1)
$styleSource = $total.Cells.Item(1, 2)

# ---------------------------------------------------------------------------
# 1. Update the "总计" (summary) sheet: insert a new row 2 for 2022-Q4 and
#    renumber the existing index column.
# ---------------------------------------------------------------------------
$total.Rows.Item(2).Insert()
$total.Range("A2:D2").ClearFormats()

$styleSource.Copy()
$total.Cells.Item(2, 1).PasteSpecial(-4122)   # xlPasteFormats

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q4"
$total.Cells.Item(2, 3).Value = 12
$total.Cells.Item(2, 4).Value = 0.73

$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(4, 1).Value = 2
$total.Cells.Item(5, 1).Value = 3
$total.Cells.Item(6, 1).Value = 4

# ---------------------------------------------------------------------------
# 2. Insert the new "2022-Q4" worksheet right after "总计" (so it becomes the
#    second tab, same spot sheet2.xml occupies in the target workbook).
# ---------------------------------------------------------------------------
$q4 = $wb.Worksheets.Add($null, $total)
$q4.Name = "2022-Q4"

# Header row.
$q4.Cells.Item(1, 2).Value = "基金代码"
$q4.Cells.Item(1, 3).Value = "基金名称"
$q4.Cells.Item(1, 4).Value = "基金规模"
$q4.Cells.Item(1, 5).Value = "股票总仓位"
$q4.Cells.Item(1, 6).Value = "仓位占比"
$q4.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q4.Cells.Item(1, 8).Value = "仓位排名"

$styleSource.Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)   # xlPasteFormats
$styleSource.Copy()
$q4.Range("A2:A13").PasteSpecial(-4122)  # xlPasteFormats

# Columns B-G hold numeric-looking text (fund codes / percentages kept as
# strings, matching the source data) so force text formatting before writing
# them; this keeps leading zeros (fund codes) intact and stops Excel from
# auto-coercing the percentage-looking strings into numbers. Row 13's "持有
# 市值(亿元)" (column G) is the sole exception: its value is a genuine 0,
# stored as a real number (matches the other quarterly sheets, where a 0.00
# holding is written as numeric 0) — so it is left out of the text range.
$q4.Range("B2:G12").NumberFormat = "@"
$q4.Range("B13:F13").NumberFormat = "@"

$rows = @(
    @("003986", "申万菱信中证500指数优选增强A",       "24.42", "93.74", "1.63", "0.3980", 9),
    @("001637", "嘉实量化精选股票",                     "12.82", "92.31", "1.16", "0.1487", 4),
    @("360001", "光大保德信量化股票",                   "12.20", "88.99", "0.78", "0.0952", 1),
    @("007794", "申万菱信中证500指数优选增强C",         "2.90",  "93.74", "1.63", "0.0473", 9),
    @("004481", "华宝第三产业灵活配置混合A",            "0.59",  "94.36", "1.91", "0.0113", 8),
    @("011389", "国都聚成混合",                         "0.35",  "81.25", "2.99", "0.0105", 8),
    @("001641", "富国绝对收益多策略定期开放混合A",      "0.93",  "70.55", "0.71", "0.0066", 6),
    @("002020", "国都创新驱动灵活配置混合",             "0.12",  "83.47", "3.28", "0.0039", 6),
    @("011231", "光大保德信锦弘混合A",                  "1.25",  "29.26", "0.31", "0.0039", 2),
    @("009149", "富国绝对收益多策略定期开放混合C",      "0.31",  "70.55", "0.71", "0.0022", 6),
    @("011232", "光大保德信锦弘混合C",                  "0.71",  "29.26", "0.31", "0.0022", 2),
    @("012798", "华宝第三产业灵活配置混合C",            "0.00",  "94.36", "1.91", 0,        8)
)

$r = 2
foreach ($item in $rows) {
    $q4.Cells.Item($r, 1).Value = $r - 2
    $q4.Cells.Item($r, 2).Value = $item[0]
    $q4.Cells.Item($r, 3).Value = $item[1]
    $q4.Cells.Item($r, 4).Value = $item[2]
    $q4.Cells.Item($r, 5).Value = $item[3]
    $q4.Cells.Item($r, 6).Value = $item[4]
    $q4.Cells.Item($r, 7).Value = $item[5]
    $q4.Cells.Item($r, 8).Value = $item[6]
    $r = $r + 1
}
